$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with P1=14, Q1=15, copying style/format from O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$excel.CutCopyMode = 0

# Row 2
$ws.Range("B2").Value = 24.95059154083955
$ws.Range("C2").Value = 20.06877886674375
$ws.Range("D2").Value = 9.611003004978826
$ws.Range("E2").Value = 30.04024143633222
$ws.Range("F2").Value = 62.62184248115091
$ws.Range("G2").Value = 2.047411979723615
$ws.Range("H2").ClearContents()
$ws.Range("I2").Value = 3.968408593914997
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0

# Row 3
$ws.Range("B3").Value = 23.16697905829687
$ws.Range("C3").Value = 18.55492189886247
$ws.Range("D3").Value = 9.090117789138334
$ws.Range("E3").Value = 27.73463633535916
$ws.Range("F3").Value = 58.55250235741687
$ws.Range("G3").Value = 2.061114558592859
$ws.Range("H3").ClearContents()
$ws.Range("I3").Value = 3.469401822719058
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0

# Row 4
$ws.Range("B4").Value = 22.01767150050155
$ws.Range("C4").Value = 17.59605588522804
$ws.Range("D4").Value = 8.756510883322148
$ws.Range("E4").Value = 26.26761867079037
$ws.Range("F4").Value = 55.94403049702101
$ws.Range("G4").Value = 2.069629128883681
$ws.Range("H4").ClearContents()
$ws.Range("I4").Value = 3.159878317903202
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0

# Row 5
$ws.Range("B5").Value = 21.53465770339592
$ws.Range("C5").Value = 17.21176187561178
$ws.Range("D5").Value = 8.59934117985527
$ws.Range("E5").Value = 25.65562908989511
$ws.Range("F5").Value = 54.7763154602218
$ws.Range("G5").Value = 2.073165267945906
$ws.Range("H5").ClearContents()
$ws.Range("I5").Value = 3.03214591473613
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0

# Row 6
$ws.Range("B6").Value = 21.45271866951602
$ws.Range("C6").Value = 17.16490513056369
$ws.Range("D6").Value = 8.551361115453581
$ws.Range("E6").Value = 25.55258073977514
$ws.Range("F6").Value = 54.48676992597131
$ws.Range("G6").Value = 2.073797108809492
$ws.Range("H6").ClearContents()
$ws.Range("I6").Value = 3.010309898759368
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0

# Row 7
$ws.Range("B7").Value = 22.00873937241082
$ws.Range("C7").Value = 17.63702284210588
$ws.Range("D7").Value = 8.696038570198752
$ws.Range("E7").Value = 26.25776914509509
$ws.Range("F7").Value = 55.67681204830991
$ws.Range("G7").Value = 2.069793429437988
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 3.156094852201096
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0

# Row 8
$ws.Range("B8").Value = 24.34308979305988
$ws.Range("C8").Value = 19.60740477918525
$ws.Range("D8").Value = 9.360647453540675
$ws.Range("E8").Value = 29.25289858800457
$ws.Range("F8").Value = 60.92969100286781
$ws.Range("G8").Value = 2.052276056651788
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 3.792948885150529
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0

# Row 9
$ws.Range("B9").Value = 28.51220527708041
$ws.Range("C9").Value = 23.19012691006922
$ws.Range("D9").Value = 10.6552598898129
$ws.Range("E9").Value = 34.76331950442783
$ws.Range("F9").Value = 70.79289306852576
$ws.Range("G9").Value = 2.018197841580306
$ws.Range("H9").ClearContents()
$ws.Range("I9").Value = 5.034186733815996
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0

# Row 10
$ws.Range("B10").Value = 31.34823089041888
$ws.Range("C10").Value = 25.74797107448514
$ws.Range("D10").Value = 11.5238667744354
$ws.Range("E10").Value = 38.65643010765509
$ws.Range("F10").Value = 77.47453495590825
$ws.Range("G10").Value = 1.99298643577381
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = 5.95454836441022
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0

# Row 11
$ws.Range("B11").Value = 32.59076895400146
$ws.Range("C11").Value = 26.94837762867594
$ws.Range("D11").Value = 11.83146400046651
$ws.Range("E11").Value = 40.4086773985107
$ws.Range("F11").Value = 80.1184638196125
$ws.Range("G11").Value = 1.981530200075788
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = 6.374139711305668
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# Row 12
$ws.Range("B12").Value = 33.05913011279707
$ws.Range("C12").Value = 27.37652983048417
$ws.Range("D12").Value = 11.98928618133663
$ws.Range("E12").Value = 41.07521670544514
$ws.Range("F12").Value = 81.28302279646911
$ws.Range("G12").Value = 1.977034313512059
$ws.Range("H12").ClearContents()
$ws.Range("I12").Value = 6.538046398285531
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0

# Row 13
$ws.Range("B13").Value = 32.95912081462323
$ws.Range("C13").Value = 27.27929406823836
$ws.Range("D13").Value = 11.96372203230739
$ws.Range("E13").Value = 40.93220356058239
$ws.Range("F13").Value = 81.06659018957856
$ws.Range("G13").Value = 1.977981847756414
$ws.Range("H13").ClearContents()
$ws.Range("I13").Value = 6.50333703022212
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0

# Row 14
$ws.Range("B14").Value = 32.62966990963741
$ws.Range("C14").Value = 26.98143693898982
$ws.Range("D14").Value = 11.84804571273014
$ws.Range("E14").Value = 40.46373884097984
$ws.Range("F14").Value = 80.22897490796537
$ws.Range("G14").Value = 1.981151614672319
$ws.Range("H14").ClearContents()
$ws.Range("I14").Value = 6.387866123358314
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Row 15
$ws.Range("B15").Value = 32.42598993564806
$ws.Range("C15").Value = 26.80917259842181
$ws.Range("D15").Value = 11.76040012624363
$ws.Range("E15").Value = 40.17579522546377
$ws.Range("F15").Value = 79.64702712033608
$ws.Range("G15").Value = 1.983131446072713
$ws.Range("H15").ClearContents()
$ws.Range("I15").Value = 6.316136037252698
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0

# Row 16
$ws.Range("B16").Value = 31.2562322021707
$ws.Range("C16").Value = 25.75837039075594
$ws.Range("D16").Value = 11.36122583403474
$ws.Range("E16").Value = 38.53303574413617
$ws.Range("F16").Value = 76.7195148798915
$ws.Range("G16").Value = 1.994101618075254
$ws.Range("H16").ClearContents()
$ws.Range("I16").Value = 5.916645168127635
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0

# Row 17
$ws.Range("B17").Value = 30.52860603931145
$ws.Range("C17").Value = 25.11260824657765
$ws.Range("D17").Value = 11.11417586347915
$ws.Range("E17").Value = 37.52343715219463
$ws.Range("F17").Value = 74.90413621496364
$ws.Range("G17").Value = 2.000765994323129
$ws.Range("H17").ClearContents()
$ws.Range("I17").Value = 5.674159198722121
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0

# Row 18
$ws.Range("B18").Value = 30.10904335535481
$ws.Range("C18").Value = 24.71196002021476
$ws.Range("D18").Value = 11.01508141054049
$ws.Range("E18").Value = 36.94382275984903
$ws.Range("F18").Value = 74.03220724564041
$ws.Range("G18").Value = 2.004470601181214
$ws.Range("H18").ClearContents()
$ws.Range("I18").Value = 5.537968254988689
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0

# Row 19
$ws.Range("B19").Value = 29.9640665381946
$ws.Range("C19").Value = 24.59799224292251
$ws.Range("D19").Value = 10.94755542456199
$ws.Range("E19").Value = 36.74539478601159
$ws.Range("F19").Value = 73.59695846622095
$ws.Range("G19").Value = 2.005806854568524
$ws.Range("H19").ClearContents()
$ws.Range("I19").Value = 5.48997992077224
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Row 20
$ws.Range("B20").Value = 30.60669339439992
$ws.Range("C20").Value = 25.17923435505836
$ws.Range("D20").Value = 11.14398425005814
$ws.Range("E20").Value = 37.63123322179111
$ws.Range("F20").Value = 75.11217609992865
$ws.Range("G20").Value = 2.000049965955155
$ws.Range("H20").ClearContents()
$ws.Range("I20").Value = 5.700087670178486
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0

# Row 21
$ws.Range("B21").Value = 32.7229971799303
$ws.Range("C21").Value = 27.09725400072911
$ws.Range("D21").Value = 11.83493828099227
$ws.Range("E21").Value = 40.59803099903036
$ws.Range("F21").Value = 80.28409172105951
$ws.Range("G21").Value = 1.980350605922608
$ws.Range("H21").ClearContents()
$ws.Range("I21").Value = 6.418130598133569
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0

# Row 22
$ws.Range("B22").Value = 34.08129923487518
$ws.Range("C22").Value = 28.31632421249627
$ws.Range("D22").Value = 12.34063268253754
$ws.Range("E22").Value = 42.5450686828101
$ws.Range("F22").Value = 83.85294153489181
$ws.Range("G22").Value = 1.967025420383101
$ws.Range("H22").ClearContents()
$ws.Range("I22").Value = 6.902902282270286
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0

# Row 23
$ws.Range("B23").Value = 33.36312412851472
$ws.Range("C23").Value = 27.63157980812691
$ws.Range("D23").Value = 12.12759591186331
$ws.Range("E23").Value = 41.50872955452025
$ws.Range("F23").Value = 82.18101189308214
$ws.Range("G23").Value = 1.974013434103063
$ws.Range("H23").ClearContents()
$ws.Range("I23").Value = 6.647584658337502
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0

# Row 24
$ws.Range("B24").Value = 30.57761457115394
$ws.Range("C24").Value = 25.09016676035566
$ws.Range("D24").Value = 11.22277727874402
$ws.Range("E24").Value = 37.58783342207083
$ws.Range("F24").Value = 75.39577328773804
$ws.Range("G24").Value = 2.000139213991619
$ws.Range("H24").ClearContents()
$ws.Range("I24").Value = 5.694405671848928
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0

# Row 25
$ws.Range("B25").Value = 27.42068120127078
$ws.Range("C25").Value = 22.30475227935743
$ws.Range("D25").Value = 10.22152264775182
$ws.Range("E25").Value = 33.30120591419096
$ws.Range("F25").Value = 67.81943386746248
$ws.Range("G25").Value = 2.027633955245948
$ws.Range("H25").ClearContents()
$ws.Range("I25").Value = 4.692282093865777
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
